$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 826.5501353333334
$ws.Range("H2").Value = 2479.650406
$ws.Range("I2").Value = 0.2851537905755522
$ws.Range("J2").Value = 0.2851537905755522
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 31.82741333333333
$ws.Range("N2").Value = 95.48223999999999
$ws.Range("O2").Value = 0.114390792932228
$ws.Range("P2").Value = 0.114390792932228
$ws.Range("Q2").Value = 26306.9527979766
$ws.Range("R2").Value = 236762.5751817894
$ws.Range("S2").Value = 0.0326189682115679
$ws.Range("T2").Value = 0.03261896821156791

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 826.5501353333334
$ws.Range("H3").Value = 2479.650406
$ws.Range("I3").Value = 0.2851537905755522
$ws.Range("J3").Value = 0.2851537905755522
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 85.46317833333335
$ws.Range("N3").Value = 256.389535
$ws.Range("O3").Value = 0.307162904935779
$ws.Range("P3").Value = 0.307162904935779
$ws.Range("Q3").Value = 70639.60161743348
$ws.Range("R3").Value = 635756.4145569013
$ws.Range("S3").Value = 0.08758866666663538
$ws.Range("T3").Value = 0.08758866666663538

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 826.5501353333334
$ws.Range("H4").Value = 2479.650406
$ws.Range("I4").Value = 0.2851537905755522
$ws.Range("J4").Value = 0.2851537905755522
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 122.2478306666667
$ws.Range("N4").Value = 366.743492
$ws.Range("O4").Value = 0.4393704929064738
$ws.Range("P4").Value = 0.4393704929064738
$ws.Range("Q4").Value = 101043.9609817398
$ws.Range("R4").Value = 909395.6488356579
$ws.Range("S4").Value = 0.1252881615193298
$ws.Range("T4").Value = 0.1252881615193298

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 826.5501353333334
$ws.Range("H5").Value = 2479.650406
$ws.Range("I5").Value = 0.2851537905755522
$ws.Range("J5").Value = 0.2851537905755522
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 38.69562533333333
$ws.Range("N5").Value = 116.086876
$ws.Range("O5").Value = 0.1390758092255191
$ws.Range("P5").Value = 0.1390758092255191
$ws.Range("Q5").Value = 31983.87435607463
$ws.Range("R5").Value = 287854.8692046716
$ws.Range("S5").Value = 0.03965799417801913
$ws.Range("T5").Value = 0.03965799417801913

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 977.7211913333334
$ws.Range("H6").Value = 2933.163574
$ws.Range("I6").Value = 0.3373067064132887
$ws.Range("J6").Value = 0.3373067064132887
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 31.82741333333333
$ws.Range("N6").Value = 95.48223999999999
$ws.Range("O6").Value = 0.114390792932228
$ws.Range("P6").Value = 0.114390792932228
$ws.Range("Q6").Value = 31118.33648132508
$ws.Range("R6").Value = 280065.0283319257
$ws.Range("S6").Value = 0.03858478160797433
$ws.Range("T6").Value = 0.03858478160797434

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 977.7211913333334
$ws.Range("H7").Value = 2933.163574
$ws.Range("I7").Value = 0.3373067064132887
$ws.Range("J7").Value = 0.3373067064132887
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 85.46317833333335
$ws.Range("N7").Value = 256.389535
$ws.Range("O7").Value = 0.307162904935779
$ws.Range("P7").Value = 0.307162904935779
$ws.Range("Q7").Value = 83559.1605351998
$ws.Range("R7").Value = 752032.4448167982
$ws.Range("S7").Value = 0.1036081077962257
$ws.Range("T7").Value = 0.1036081077962257

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 977.7211913333334
$ws.Range("H8").Value = 2933.163574
$ws.Range("I8").Value = 0.3373067064132887
$ws.Range("J8").Value = 0.3373067064132887
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 122.2478306666667
$ws.Range("N8").Value = 366.743492
$ws.Range("O8").Value = 0.4393704929064738
$ws.Range("P8").Value = 0.4393704929064738
$ws.Range("Q8").Value = 119524.2946373289
$ws.Range("R8").Value = 1075718.651735961
$ws.Range("S8").Value = 0.1482026138574659
$ws.Range("T8").Value = 0.1482026138574659

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 977.7211913333334
$ws.Range("H9").Value = 2933.163574
$ws.Range("I9").Value = 0.3373067064132887
$ws.Range("J9").Value = 0.3373067064132887
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 38.69562533333333
$ws.Range("N9").Value = 116.086876
$ws.Range("O9").Value = 0.1390758092255191
$ws.Range("P9").Value = 0.1390758092255191
$ws.Range("Q9").Value = 37833.53290029498
$ws.Range("R9").Value = 340501.7961026548
$ws.Range("S9").Value = 0.04691120315162273
$ws.Range("T9").Value = 0.04691120315162273

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 577.0637716666666
$ws.Range("H10").Value = 1731.191315
$ws.Range("I10").Value = 0.1990828079995583
$ws.Range("J10").Value = 0.1990828079995583
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 31.82741333333333
$ws.Range("N10").Value = 95.48223999999999
$ws.Range("O10").Value = 0.114390792932228
$ws.Range("P10").Value = 0.114390792932228
$ws.Range("Q10").Value = 18366.44718052728
$ws.Range("R10").Value = 165298.0246247456
$ws.Range("S10").Value = 0.02277324026624397
$ws.Range("T10").Value = 0.02277324026624398

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 577.0637716666666
$ws.Range("H11").Value = 1731.191315
$ws.Range("I11").Value = 0.1990828079995583
$ws.Range("J11").Value = 0.1990828079995583
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 85.46317833333335
$ws.Range("N11").Value = 256.389535
$ws.Range("O11").Value = 0.307162904935779
$ws.Range("P11").Value = 0.307162904935779
$ws.Range("Q11").Value = 49317.70402765428
$ws.Range("R11").Value = 443859.3362488886
$ws.Range("S11").Value = 0.06115085362791627
$ws.Range("T11").Value = 0.06115085362791627

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 577.0637716666666
$ws.Range("H12").Value = 1731.191315
$ws.Range("I12").Value = 0.1990828079995583
$ws.Range("J12").Value = 0.1990828079995583
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 122.2478306666667
$ws.Range("N12").Value = 366.743492
$ws.Range("O12").Value = 0.4393704929064738
$ws.Range("P12").Value = 0.4393704929064738
$ws.Range("Q12").Value = 70544.79424257466
$ws.Range("R12").Value = 634903.148183172
$ws.Range("S12").Value = 0.0874711114799708
$ws.Range("T12").Value = 0.08747111147997082

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 577.0637716666666
$ws.Range("H13").Value = 1731.191315
$ws.Range("I13").Value = 0.1990828079995583
$ws.Range("J13").Value = 0.1990828079995583
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 38.69562533333333
$ws.Range("N13").Value = 116.086876
$ws.Range("O13").Value = 0.1390758092255191
$ws.Range("P13").Value = 0.1390758092255191
$ws.Range("Q13").Value = 22329.84350185355
$ws.Range("R13").Value = 200968.5915166819
$ws.Range("S13").Value = 0.02768760262542722
$ws.Range("T13").Value = 0.02768760262542722

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 517.2766776666667
$ws.Range("H14").Value = 1551.830033
$ws.Range("I14").Value = 0.1784566950116009
$ws.Range("J14").Value = 0.1784566950116009
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 31.82741333333333
$ws.Range("N14").Value = 95.48223999999999
$ws.Range("O14").Value = 0.114390792932228
$ws.Range("P14").Value = 0.114390792932228
$ws.Range("Q14").Value = 16463.57862779043
$ws.Range("R14").Value = 148172.2076501139
$ws.Range("S14").Value = 0.0204138028464418
$ws.Range("T14").Value = 0.02041380284644181

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 517.2766776666667
$ws.Range("H15").Value = 1551.830033
$ws.Range("I15").Value = 0.1784566950116009
$ws.Range("J15").Value = 0.1784566950116009
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 85.46317833333335
$ws.Range("N15").Value = 256.389535
$ws.Range("O15").Value = 0.307162904935779
$ws.Range("P15").Value = 0.307162904935779
$ws.Range("Q15").Value = 44208.10895110053
$ws.Range("R15").Value = 397872.9805599047
$ws.Range("S15").Value = 0.05481527684500166
$ws.Range("T15").Value = 0.05481527684500166

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 517.2766776666667
$ws.Range("H16").Value = 1551.830033
$ws.Range("I16").Value = 0.1784566950116009
$ws.Range("J16").Value = 0.1784566950116009
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 122.2478306666667
$ws.Range("N16").Value = 366.743492
$ws.Range("O16").Value = 0.4393704929064738
$ws.Range("P16").Value = 0.4393704929064738
$ws.Range("Q16").Value = 63235.95169921059
$ws.Range("R16").Value = 569123.5652928953
$ws.Range("S16").Value = 0.07840860604970733
$ws.Range("T16").Value = 0.07840860604970734

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 517.2766776666667
$ws.Range("H17").Value = 1551.830033
$ws.Range("I17").Value = 0.1784566950116009
$ws.Range("J17").Value = 0.1784566950116009
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 38.69562533333333
$ws.Range("N17").Value = 116.086876
$ws.Range("O17").Value = 0.1390758092255191
$ws.Range("P17").Value = 0.1390758092255191
$ws.Range("Q17").Value = 20016.34451266077
$ws.Range("R17").Value = 180147.1006139469
$ws.Range("S17").Value = 0.02481900927045005
$ws.Range("T17").Value = 0.02481900927045005
